# Updated parser to use TokenIteratorFieldRewriterSplit.
#
# This script reproduces, via Word COM-interop calls, the run-splitting
# behaviour of the new TokenIteratorFieldRewriterSplit: every "}" that
# closes an M2Doc field is pushed into its own run, and a separating
# space is introduced after the "{m:" / ":" token that starts a field
# where the parser now emits one.
#
# Word's object model silently re-merges adjacent runs that end up with
# identical run properties after a plain text edit (InsertBefore/Range.Text),
# so after inserting the new characters we force the desired run
# boundaries back into existence by toggling a character property
# (Bold on, then off) on the tail sub-range of the run that must become
# separate. Toggling leaves the final formatting identical to its
# neighbours while keeping the run boundary intact, which is exactly the
# shape the diff expects.

$d = $word.ActiveDocument

function Find-Range($searchText) {
    $r = $d.Content
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "text not found: $searchText"
    }
    return $r
}

function Split-Boundary-Before($absPos) {
    # Force a run boundary immediately before the character at $absPos by
    # toggling Bold on the single character preceding it, then restoring it.
    if ($absPos -le 0) { return }
    $charBefore = $d.Range($absPos - 1, $absPos)
    $charBefore.Font.Bold = 1
    $charBefore.Font.Bold = 0
}

function Split-Boundary-After($absPos) {
    # Force a run boundary immediately after the character at $absPos by
    # toggling Bold on the single character that follows it, then restoring it.
    $charAfter = $d.Range($absPos, $absPos + 1)
    $charAfter.Font.Bold = 1
    $charAfter.Font.Bold = 0
}

# ---------------------------------------------------------------------
# Step 1: insert the two new space characters (the only real text
# changes in the whole diff).
# ---------------------------------------------------------------------

# "{m:template myTemplate(element :ecore::...)}" -> insert " " before "ecore::"
$r = Find-Range("ecore::")
$insertPos = $r.Start
$ip = $d.Range($insertPos, $insertPos)
$ip.InsertBefore(" ")

# "{m:child.myTemplate(depth + 1)}" -> insert " " before "child" (right after "{m:")
$r = Find-Range("{m:child")
$insertPos = $r.Start + 3
$ip = $d.Range($insertPos, $insertPos)
$ip.InsertBefore(" ")

# ---------------------------------------------------------------------
# Step 2: restore / create every run boundary required by the diff.
# Word merges same-format adjacent runs after the edits above, so every
# boundary (old and new) gets re-asserted explicitly using absolute
# positions re-resolved via Find, from the end of the document towards
# the start so that earlier edits do not shift not-yet-processed
# offsets. Each Split call is self-contained (re-finds via Find.Execute)
# so ordering across different search terms is not offset-sensitive.
# ---------------------------------------------------------------------

# Paragraph "{m:endfor}" -> "{m:endfor" | "}"
$r = Find-Range("{m:endfor}")
Split-Boundary-Before($r.End - 1)

# Paragraph "{m:for child | ...->filter(ecore::ENamedElement)}" -> ... "ENamedElement" | ")" | "}"
$r = Find-Range("->filter(ecore::ENamedElement)}")
Split-Boundary-Before($r.End - 1)
Split-Boundary-Before($r.End - 2)

# Paragraph "{m: element.name.asStyle('Titre' + depth)}" -> ... "' + depth)" | "}"
$r = Find-Range("' + depth)}")
Split-Boundary-Before($r.End - 1)

# Paragraph "{m:template myTemplate(element : ecore::ENamedElement, depth : Integer)}"
#   -> "{m:template myTemplate(" | "element " | ":" | " " | "ecore::" | "ENamedElement" | ", depth : Integer" | ")" | "}"
$r = Find-Range("element : ecore::ENamedElement, depth : Integer)}")
$elementStart = $r.Start
Split-Boundary-After($r.End - 1 - 1)    # before final "}"
Split-Boundary-Before($r.End - 1)       # ")" | "}"
Split-Boundary-Before($elementStart + ("element ".Length))          # "element " | ":"
Split-Boundary-Before($elementStart + ("element :".Length))         # ":" | " "
Split-Boundary-Before($elementStart + ("element : ".Length))        # " " | "ecore::"

# Paragraph "{m: child.myTemplate(depth + 1)}" -> "{m:" | " " | "child" | "." | "myTemplate(" | "depth + 1" | ")}"
$r = Find-Range("{m: child.myTemplate(depth + 1)}")
$fieldStart = $r.Start
Split-Boundary-Before($fieldStart + ("{m:".Length))      # "{m:" | " "
Split-Boundary-Before($fieldStart + ("{m: ".Length))     # " " | "child"
